$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold plain text in this sheet, not
# numbers (e.g. "37.557.47", "  +0.74%  "). Several column D strings look like
# valid numeric literals (e.g. "263.70", "55.99"), so assigning them directly
# would make Excel silently reinterpret the cell as a number. To keep them as
# text we temporarily force Text format before writing the value, then restore
# the default "Normal" style afterwards so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.557.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.016.67'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '263.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.618'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.14%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.99'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.386'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("E10").Value = '  -3.66%  '
$ws.Range("E11").Value = '  -2.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.35'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.309.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.806'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.33%  '
$ws.Range("E16").Value = '  -3.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.004.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.504.39'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0843'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.97'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.26%  '
$ws.Range("E23").Value = '  +7.89%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("E27").Value = '  -5.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.128'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -11.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.33'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.50%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0651'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  +1.29%  '
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.18'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.13%  '
$ws.Range("E40").Value = '  +4.86%  '
$ws.Range("E41").Value = '  +2.97%  '
$ws.Range("E42").Value = '  -4.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0214'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.396.63'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '89.95'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.67'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.48%  '
$ws.Range("E47").Value = '  -1.89%  '
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.202.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.95'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.11%  '
